# Add a new data row (row 22) to the active sheet, mirroring the same
# fixed column layout used by the existing rows (e.g. row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22

$ws.Cells.Item($row, 1).Value  = 131271025          # A  Id
$ws.Cells.Item($row, 2).Value  = 57881              # B  Taxonsorteringsordning
$ws.Cells.Item($row, 4).Value  = "NT"                # D  Rödlistade
$ws.Cells.Item($row, 5).Value  = 100049             # E  TaxonId
$ws.Cells.Item($row, 6).Value  = "Spillkråka"        # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Dryocopus martius" # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(Linnaeus, 1758)"  # H  Auktor

# Columns that exist in the source export template but carry no data for
# this observation are still written out as (explicitly typed) empty
# text cells rather than left completely blank. A lone quote-prefix
# produces an empty-string text cell; the style is then reset so no
# quote-prefix formatting lingers on the cell.
$emptyTextCols = 9, 11, 12, 46, 51   # I, K, L, AT, AY
foreach ($col in $emptyTextCols) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'"
    $c.Style = "Normal"
}

$ws.Cells.Item($row, 13).Value = "spel/sång"         # M  Aktivitet
$ws.Cells.Item($row, 14).Value = "observerad"        # N  Metod

$ws.Cells.Item($row, 16).Value = "Vid Lillebo, Ög"   # P  Lokalnamn
$ws.Cells.Item($row, 17).Value = 567427             # Q  Ost
$ws.Cells.Item($row, 18).Value = 6510021            # R  Nord
$ws.Cells.Item($row, 19).Value = 25                 # S  Noggrannhet
$ws.Cells.Item($row, 20).Value = "Östergötland"      # T  Län
$ws.Cells.Item($row, 21).Value = "Norrköping"        # U  Kommun
$ws.Cells.Item($row, 22).Value = "Östergötland"      # V  Provins
$ws.Cells.Item($row, 23).Value = "Simonstorp"        # W  Socken

# Y (25) / AA (27) hold plain-text dates ("2026-02-23"), not Excel date
# serials, so force a text number format before assigning, then restore
# the Normal style so no stray formatting is left on the cell.
$dateCols = 25, 27
foreach ($col in $dateCols) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = "2026-02-23"
    $c.Style = "Normal"
}

$ws.Cells.Item($row, 30).Value = $false   # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false   # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false   # AG Ospontan

$ws.Cells.Item($row, 49).Value = "Anette Källman"  # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Anette Källman"  # AX Observatörer
